# Weekly update: insert a new price record for "Hortaliza" (Choclo) as row 107,
# pushing the existing rows 107-202 down to 108-203.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 107 (this shifts rows 107:202 -> 108:203 and extends the
# used range from A1:R202 to A1:R203, carrying formatting down with it).
$ws.Rows("107:107").Insert()

# Populate the newly inserted row 107 with the new weekly record.
$ws.Cells.Item(107, 1).Value  = 2
$ws.Cells.Item(107, 2).Value  = "Comercializadora del Agro de Limarí"
$ws.Cells.Item(107, 3).Value  = "Coquimbo"
$ws.Cells.Item(107, 4).Value  = 45225
$ws.Cells.Item(107, 5).Value  = 4
$ws.Cells.Item(107, 6).Value  = 100112024
$ws.Cells.Item(107, 7).Value  = "Choclo"
$ws.Cells.Item(107, 8).Value  = "Dulce o Americano"
$ws.Cells.Item(107, 9).Value  = "Primera"
$ws.Cells.Item(107, 10).Value = 500
$ws.Cells.Item(107, 11).Value = 26000
$ws.Cells.Item(107, 12).Value = 28000
$ws.Cells.Item(107, 13).Value = 27000
$ws.Cells.Item(107, 14).Value = "$/malla 70 unidades"
$ws.Cells.Item(107, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(107, 16).Value = 386
$ws.Cells.Item(107, 17).Value = 70
$ws.Cells.Item(107, 18).Value = "Hortaliza"
